# Cambios para el excel
# Adds 4 new cells (M1:P1) to row 1, extends the used range to A1:P4, and
# moves the sheet's view/selection to reflect the new columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values at the end of row 1 (columns M, N, O, P).
$ws.Range("M1").Value = 7980
$ws.Range("N1").Value = "calle 01"
$ws.Range("O1").Value = 23
$ws.Range("P1").Value = 34

# Activate the sheet and move the view so the new columns are visible,
# then leave the selection on O3 (matches the target selection state).
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("O3").Select()
